# Generate Report for Handback
# Update the timestamp values recorded on the handback status sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview!G2 - Latest HO Xliff Generate Date
$wsOverview.Range("G2").Value = "2016-08-19 15:12:34"

# zh-cn!H2 - Correspond Handoff Datetime
$wsZhCn.Range("H2").Value = "2016-08-19 15:12:29"

# zh-cn!K2 - Correspond Handback DateTime
$wsZhCn.Range("K2").Value = "2016-08-19 15:12:46"

# de-de!K2 - Correspond Handback DateTime
$wsDeDe.Range("K2").Value = "2016-08-19 15:12:52"
